# Update column F (dSF) values on the active worksheet to reflect
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = -3
    4  = 5
    5  = 5
    6  = 4
    7  = 2
    8  = 6
    9  = 3
    10 = 3
    12 = 9
    13 = 3
    14 = -3
    15 = -2
    16 = 2
    18 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
